$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.243.39'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '2.613.42'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.598'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.587'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.74%  '
$ws.Range("E11").Value = '  +2.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("E13").Value = '  +2.43%  '
$ws.Range("D14").Value = '3.007.35'
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("E15").Value = '  +2.01%  '
$ws.Range("D16").Value = '2.606.39'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.921'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.93'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").Value = '46.448.13'
$ws.Range("E19").Value = '  +0.81%  '
$ws.Range("E20").Value = '  +2.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '290.65'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +15.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.00%  '
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.32%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  +1.45%  '
$ws.Range("E30").Value = '  +5.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '39.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.40%  '
$ws.Range("E34").Value = '  -1.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '157.77'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.69%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.22'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.85%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0841'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.49%  '
$ws.Range("E39").Value = '  +7.19%  '
$ws.Range("E40").Value = '  +2.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0334'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.60'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.29%  '
$ws.Range("E44").Value = '  -2.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +13.02%  '
$ws.Range("D46").Value = '2.115.42'
$ws.Range("E46").Value = '  +3.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.51%  '
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.51'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '109.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("E51").Value = '  +2.90%  '
